# Update financial figures on the GLAE worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLAE")

# Cells that become "NA" (text) values
$ws.Range("D8").Value = "NA"
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D18").Value = "NA"
$ws.Range("D20").Value = "NA"
$ws.Range("J21").Value = "NA"
$ws.Range("D32").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"

# Cell that becomes numeric (was "NA")
$ws.Range("D22").Value = 0

# Cells with updated numeric values
$ws.Range("D14").Value = 8200
$ws.Range("D17").Value = 8700
$ws.Range("D21").Value = -3900
$ws.Range("D23").Value = -8000
$ws.Range("D24").Value = -5700
$ws.Range("D26").Value = -2300
$ws.Range("D27").Value = 7900
$ws.Range("D29").Value = -111000
$ws.Range("D33").Value = -103100
$ws.Range("D35").Value = -103100
$ws.Range("D81").Value = -103100
$ws.Range("I91").Value = -10200
